$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.689.70"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "3.258.98"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.259.82"
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.503"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").Value = "3.785.96"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").Value = "66.734.87"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "3.254.64"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.113"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "507.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.755"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.135"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +53.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +21.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("D39").Value = "0.0₃0784"
$ws.Range("E39").Value = "  +16.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "495.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0427"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.295"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.09%  "
$ws.Range("D46").Value = "2.995.99"
$ws.Range("E46").Value = "  +6.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.43%  "
$ws.Range("E48").Value = "  +6.40%  "
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
